# NR12_temp.docx — "adj remocao hse Bruna"
#
# The HSE signatory block in the signature table is being swapped from
# Bruna Petroni Cesário / Gerente de HSE Brasil
# to
# LEONARDO SILVERIO FERREIRA / Técnico(a) de Segurança do Trabalho
#
# Both strings are unique in the document body (the name is a bold
# paragraph, the title is the plain paragraph right below it inside the
# same table cell), so a simple formatted Find/Replace on each locates
# the correct run without disturbing anything else (table layout, the
# other column with "RENNNNNN" / role / CPF, etc.).
#
# Find.Execute Replace:=2 (wdReplaceOne) reuses the formatting already
# present on the matched run, so the existing bold / Arial run
# properties on these two runs are left intact automatically.

$d = $word.ActiveDocument

# Build the accented strings from code points so the script is robust
# regardless of how this file's bytes get transported.
$aAcute     = [char]0x00E1   # á
$eAcute     = [char]0x00E9   # é
$cCedilla   = [char]0x00E7   # ç

$oldName = "Bruna Petroni Ces" + $aAcute + "rio"
$newName = "LEONARDO SILVERIO FERREIRA"

$oldTitle = "Gerente de HSE Brasil"
$newTitle = "T" + $eAcute + "cnico(a) de Seguran" + $cCedilla + "a do Trabalho"

$foundName = $d.Content.Find.Execute(
    $oldName, $true, $false, $false, $false, $false,
    $true, 1, $false, $newName, 2)

$foundTitle = $d.Content.Find.Execute(
    $oldTitle, $true, $false, $false, $false, $false,
    $true, 1, $false, $newTitle, 2)

Write-Output ("Replaced signatory name: " + $foundName)
Write-Output ("Replaced signatory title: " + $foundTitle)
